$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-07 Friday" "2025-03-08 Saturday"

Replace-Text "428×3=1284" "868×7=6076"
Replace-Text "813×3=2439" "711×5=3555"
Replace-Text "734×4=2936" "205×8=1640"
Replace-Text "360×2=720" "646×9=5814"
Replace-Text "541×8=4328" "678×3=2034"
Replace-Text "821×4=3284" "493×6=2958"
Replace-Text "542×3=1626" "723×8=5784"
Replace-Text "847×7=5929" "347×4=1388"
Replace-Text "743×6=4458" "596×7=4172"
Replace-Text "836×8=6688" "171×8=1368"
Replace-Text "983×8=7864" "697×5=3485"
Replace-Text "277×2=554" "204×9=1836"
Replace-Text "516×5=2580" "972×8=7776"
Replace-Text "178×8=1424" "599×6=3594"
Replace-Text "485×7=3395" "235×4=940"
Replace-Text "393×3=1179" "466×6=2796"
Replace-Text "929×3=2787" "602×7=4214"
Replace-Text "410×2=820" "828×8=6624"
Replace-Text "332×9=2988" "609×3=1827"
Replace-Text "709×3=2127" "574×2=1148"
Replace-Text "603×2=1206" "181×6=1086"
Replace-Text "163×4=652" "601×6=3606"
Replace-Text "962×5=4810" "911×4=3644"
Replace-Text "527×3=1581" "475×4=1900"
Replace-Text "383×8=3064" "581×5=2905"
